$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E46").Value = 91.37
$ws.Range("E49").Value = 98.08
$ws.Range("E50").Value = 98.39
$ws.Range("E51").Value = 100.3
$ws.Range("E52").Value = 99.89
$ws.Range("E53").Value = 105.4
$ws.Range("E54").Value = 97.68000000000001
$ws.Range("E55").Value = 100.72
$ws.Range("E56").Value = 100.04
$ws.Range("E57").Value = 103.51
$ws.Range("E58").Value = 98.75
$ws.Range("E62").Value = 97.08
$ws.Range("C64").Value = 98.84
$ws.Range("E64").Value = 100.78
$ws.Range("E65").Value = 103.43
$ws.Range("C66").Value = 100.78
$ws.Range("C68").Value = 98.39
$ws.Range("E69").Value = 99.84999999999999
$ws.Range("C70").Value = 94.73
$ws.Range("E70").Value = 96.15000000000001
$ws.Range("C71").Value = 97.67
$ws.Range("E71").Value = 95.70999999999999
$ws.Range("C72").Value = 100.93
$ws.Range("E72").Value = 96.53
$ws.Range("C73").Value = 100.92
$ws.Range("E73").Value = 98.38
$ws.Range("C74").Value = 102.47
$ws.Range("E74").Value = 97.37
$ws.Range("C75").Value = 102.14
$ws.Range("E75").Value = 100.13
$ws.Range("C76").Value = 101.3
$ws.Range("E76").Value = 96.33
$ws.Range("C77").Value = 102.77
$ws.Range("C78").Value = 100.14
$ws.Range("E78").Value = 97.95999999999999
$ws.Range("C79").Value = 102.18
$ws.Range("E79").Value = 97.95999999999999
$ws.Range("C80").Value = 102.03
$ws.Range("E80").Value = 97.5
$ws.Range("C81").Value = 101.54
$ws.Range("E81").Value = 100.32
$ws.Range("C82").Value = 101.79
$ws.Range("E82").Value = 100.12
$ws.Range("C83").Value = 97.87
$ws.Range("E83").Value = 91.88
$ws.Range("B84").Value = 98.42
$ws.Range("C84").Value = 99.03
$ws.Range("E84").Value = 97.31
$ws.Range("C85").Value = 101.06
$ws.Range("C86").Value = 102.11
$ws.Range("E86").Value = 100.84
$ws.Range("C87").Value = 102.95
$ws.Range("E87").Value = 102.55
